$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 69753
$ws.Range("J7").Value = 69753
$ws.Range("L7").Value = 69753
$ws.Range("N7").Value = -69977
# Row 14
$ws.Range("H14").Value = 69753
$ws.Range("J14").Value = 69753
$ws.Range("L14").Value = 69753
$ws.Range("N14").Value = -70135
# Row 76
$ws.Range("H76").Value = 3258.8235
$ws.Range("I76").Value = 3314.2856
$ws.Range("K76").Value = 3314.2856
$ws.Range("M76").Value = -2999.2856
# Row 79
$ws.Range("H79").Value = 3258.8235
$ws.Range("I79").Value = 3314.2856
$ws.Range("K79").Value = 3314.2856
$ws.Range("M79").Value = -2222.2856
# Row 113
$ws.Range("H113").Value = 2650
$ws.Range("I113").Value = 2533.3333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2533.3333
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 720.6667000000002
$ws.Range("N113").Value = -9508
# Row 116
$ws.Range("H116").Value = 2699.4707
$ws.Range("I116").Value = 2661.923
$ws.Range("J116").Value = 2821.5
$ws.Range("K116").Value = 2661.923
$ws.Range("L116").Value = 2821.5
$ws.Range("M116").Value = 780.0770000000002
$ws.Range("N116").Value = -9705.5
# Row 137
$ws.Range("H137").Value = 3378.28
$ws.Range("I137").Value = 2750.1738
$ws.Range("K137").Value = 8250.5214
$ws.Range("M137").Value = -5700.5214
# Row 138
$ws.Range("H138").Value = 1963.7028
$ws.Range("I138").Value = 1526.8966
$ws.Range("J138").Value = 2245.2
$ws.Range("K138").Value = 4580.6898
$ws.Range("L138").Value = 6735.599999999999
$ws.Range("N138").Value = -17015.6
$ws.Range("M138").Value = 559.3101999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 102
$ws.Range("H102").Value = 2670.6428
$ws.Range("I102").Value = 2714.5386
$ws.Range("K102").Value = 2714.5386
$ws.Range("M102").Value = -1092.5386
# Row 128
$ws.Range("H128").Value = 34439.5
$ws.Range("J128").Value = 34439.5
$ws.Range("L128").Value = 34439.5
$ws.Range("N128").Value = -44399.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1411.2667
$ws.Range("I99").Value = 1496.5834
$ws.Range("J99").Value = 1070
$ws.Range("K99").Value = 1496.5834
$ws.Range("L99").Value = 1070
$ws.Range("M99").Value = 1.416600000000017
$ws.Range("N99").Value = -4066

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5579.483
$ws.Range("I31").Value = 1145.7222
$ws.Range("J31").Value = 12834.728
$ws.Range("K31").Value = 1145.7222
$ws.Range("L31").Value = 12834.728
$ws.Range("M31").Value = -850.7221999999999
$ws.Range("N31").Value = -13424.728
# Row 34
$ws.Range("H34").Value = 5579.483
$ws.Range("I34").Value = 1145.7222
$ws.Range("J34").Value = 12834.728
$ws.Range("K34").Value = 1145.7222
$ws.Range("L34").Value = 12834.728
$ws.Range("M34").Value = -943.7221999999999
$ws.Range("N34").Value = -13238.728
# Row 122
$ws.Range("H122").Value = 1748.05
$ws.Range("I122").Value = 1648.7826
$ws.Range("J122").Value = 1882.3529
$ws.Range("K122").Value = 4946.3478
$ws.Range("L122").Value = 5647.0587
$ws.Range("M122").Value = -2496.3478
$ws.Range("N122").Value = -10547.0587
# Row 134
$ws.Range("H134").Value = 4214.6772
$ws.Range("I134").Value = 4550.222
$ws.Range("J134").Value = 1949.75
$ws.Range("K134").Value = 13650.666
$ws.Range("L134").Value = 5849.25
$ws.Range("M134").Value = -11115.666
$ws.Range("N134").Value = -10919.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 590.6667
$ws.Range("J5").Value = 1416
$ws.Range("L5").Value = 4248
$ws.Range("N5").Value = -4472
# Row 9
$ws.Range("H9").Value = 47854.42
$ws.Range("J9").Value = 50494.5
$ws.Range("L9").Value = 151483.5
$ws.Range("N9").Value = -151931.5
# Row 122
$ws.Range("H122").Value = 8445.23
$ws.Range("J122").Value = 17464
$ws.Range("L122").Value = 157176
$ws.Range("N122").Value = -162076
# Row 132
$ws.Range("H132").Value = 1718.8695
$ws.Range("I132").Value = 1257.5
$ws.Range("J132").Value = 1964.9333
$ws.Range("K132").Value = 11317.5
$ws.Range("L132").Value = 17684.3997
$ws.Range("M132").Value = -8787.5
$ws.Range("N132").Value = -22744.3997
# Row 134
$ws.Range("H134").Value = 6361.0557
$ws.Range("I134").Value = 3306.7
$ws.Range("J134").Value = 7535.8076
$ws.Range("K134").Value = 9920.099999999999
$ws.Range("L134").Value = 22607.4228
$ws.Range("M134").Value = -4850.099999999999
$ws.Range("N134").Value = -32747.4228
# Row 135
$ws.Range("H135").Value = 590.6667
$ws.Range("J135").Value = 1416
$ws.Range("L135").Value = 12744
$ws.Range("N135").Value = -17814

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 30666.666
$ws.Range("J62").Value = 30666.666
$ws.Range("L62").Value = 30666.666
$ws.Range("N62").Value = -32038.666
# Row 65
$ws.Range("H65").Value = 30666.666
$ws.Range("J65").Value = 30666.666
$ws.Range("L65").Value = 91999.99800000001
$ws.Range("N65").Value = -98863.99800000001
# Row 80
$ws.Range("H80").Value = 50903390
$ws.Range("I80").Value = 72717850
$ws.Range("K80").Value = 72717850
$ws.Range("M80").Value = -72716852
# Row 83
$ws.Range("H83").Value = 50903390
$ws.Range("I83").Value = 72717850
$ws.Range("K83").Value = 363589250
$ws.Range("M83").Value = -363584258
# Row 120
$ws.Range("H120").Value = 32700
$ws.Range("I120").Value = 30000
$ws.Range("K120").Value = 30000
$ws.Range("M120").Value = -25162
# Row 122
$ws.Range("H122").Value = 1613.8572
$ws.Range("I122").Value = 1660.2222
$ws.Range("J122").Value = 1335.6666
$ws.Range("K122").Value = 4980.6666
$ws.Range("L122").Value = 4006.9998
$ws.Range("M122").Value = -2530.6666
$ws.Range("N122").Value = -8906.9998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1000004
$ws.Range("I40").Value = 1000004
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1000004
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -999868
$ws.Range("N40").ClearContents()
# Row 46
$ws.Range("H46").Value = 785.3
$ws.Range("I46").Value = 650.25
$ws.Range("J46").Value = 875.3333
$ws.Range("K46").Value = 650.25
$ws.Range("L46").Value = 875.3333
$ws.Range("M46").Value = -462.25
$ws.Range("N46").Value = -1251.3333
# Row 55
$ws.Range("H55").Value = 657.6667
$ws.Range("I55").Value = 478.7143
$ws.Range("J55").Value = 908.2
$ws.Range("K55").Value = 478.7143
$ws.Range("L55").Value = 908.2
$ws.Range("M55").Value = -305.7143
$ws.Range("N55").Value = -1254.2
# Row 68
$ws.Range("H68").Value = 1596
$ws.Range("I68").Value = 1490.5238
$ws.Range("K68").Value = 1490.5238
$ws.Range("M68").Value = -741.5237999999999
# Row 71
$ws.Range("H71").Value = 1596
$ws.Range("I71").Value = 1490.5238
$ws.Range("K71").Value = 7452.619
$ws.Range("M71").Value = -3708.619
# Row 82
$ws.Range("H82").Value = 166669660
$ws.Range("I82").Value = 250002000
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 250002000
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -250001639
$ws.Range("N82").Value = -5722
# Row 85
$ws.Range("H85").Value = 166669660
$ws.Range("I85").Value = 250002000
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 250002000
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -250000752
$ws.Range("N85").Value = -7496
# Row 122
$ws.Range("H122").Value = 3036.0417
$ws.Range("I122").Value = 3045.4546
$ws.Range("J122").Value = 3028.077
$ws.Range("K122").Value = 9136.363799999999
$ws.Range("L122").Value = 9084.231
$ws.Range("M122").Value = -6686.363799999999
$ws.Range("N122").Value = -13984.231
# Row 136
$ws.Range("H136").Value = 18520698
$ws.Range("I136").Value = 2228.4285
$ws.Range("J136").Value = 83335336
$ws.Range("K136").Value = 6685.2855
$ws.Range("L136").Value = 250006008
$ws.Range("M136").Value = -4135.2855
$ws.Range("N136").Value = -250011108

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4284.6
$ws.Range("J81").Value = 3812.375
$ws.Range("L81").Value = 7624.75
$ws.Range("N81").Value = -9746.75
# Row 84
$ws.Range("H84").Value = 4284.6
$ws.Range("J84").Value = 3812.375
$ws.Range("L84").Value = 38123.75
$ws.Range("N84").Value = -48731.75
# Row 122
$ws.Range("H122").Value = 2095.258
$ws.Range("I122").Value = 1841.2174
$ws.Range("K122").Value = 5523.6522
$ws.Range("M122").Value = -3073.6522
# Row 132
$ws.Range("H132").Value = 3969973
$ws.Range("I132").Value = 1584.5
$ws.Range("J132").Value = 9261158
$ws.Range("K132").Value = 4753.5
$ws.Range("L132").Value = 27783474
$ws.Range("M132").Value = -2223.5
$ws.Range("N132").Value = -27788534
